$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "NSE:AETHER"
$ws.Range("C2").Value = "NSE:ALOKINDS"
$ws.Range("D2").Value = "NSE:AXISBANK"
$ws.Range("E2").Value = ""

# Row 3
$ws.Range("B3").Value = "NSE:ASAHIINDIA"
$ws.Range("C3").Value = "NSE:ASIANTILES"

# Row 4
$ws.Range("B4").Value = "NSE:BINANIIND"
$ws.Range("C4").Value = "NSE:BEML"

# Row 5
$ws.Range("B5").Value = "NSE:CLEAN"
$ws.Range("C5").Value = "NSE:CENTRUM"

# Row 6
$ws.Range("B6").Value = "NSE:DBSTOCKBRO"
$ws.Range("C6").Value = "NSE:CERA"

# Row 7
$ws.Range("B7").Value = "NSE:DHUNINV"
$ws.Range("C7").Value = "NSE:CREATIVE"

# Row 8
$ws.Range("B8").Value = "NSE:GLOBUSSPR"
$ws.Range("C8").Value = "NSE:DBOL"

# Row 9
$ws.Range("B9").Value = "NSE:INOXWIND"
$ws.Range("C9").Value = "NSE:DHANI"

# Row 10
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = "NSE:FINOPB"

# Row 11
$ws.Range("C11").Value = "NSE:GATEWAY"

# Row 12
$ws.Range("C12").Value = "NSE:HOMEFIRST"

# Row 13
$ws.Range("C13").Value = "NSE:ICIL"

# Row 14
$ws.Range("C14").Value = "NSE:KOTHARIPET"

# Row 15
$ws.Range("C15").Value = "NSE:MANGLMCEM"

# Row 16
$ws.Range("C16").Value = "NSE:MATRIMONY"

# Row 17
$ws.Range("C17").Value = "NSE:MODISONLTD"

# Row 18
$ws.Range("C18").Value = "NSE:MURUDCERA"

# Row 19
$ws.Range("C19").Value = "NSE:PNB"

# Row 20
$ws.Range("C20").Value = "NSE:PSUBNKBEES"

# New rows 21-23
$ws.Range("A21").Value = 19
$ws.Range("C21").Value = "NSE:RML"
$ws.Range("A21").Style = $ws.Range("A20").Style

$ws.Range("A22").Value = 20
$ws.Range("C22").Value = "NSE:SAGCEM"
$ws.Range("A22").Style = $ws.Range("A20").Style

$ws.Range("A23").Value = 21
$ws.Range("C23").Value = "NSE:SAMBHAAV"
$ws.Range("A23").Style = $ws.Range("A20").Style
